$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update header figures (VALOR MORA total and Cant. Periodos count)
# ---------------------------------------------------------------------
$ws.Range("E11").Value2 = 341640
$ws.Range("F13").Value2 = 3

# ---------------------------------------------------------------------
# 2. Grow the worker detail table from 4 rows (16-19) to 6 rows (16-21).
#    Insert two new rows before the current last row (19) and copy the
#    formatting from row 18 (a normal, non-final body row) into them so
#    the borders/number-formats match the rest of the table. The former
#    last row (19, with the special bottom-border style) is pushed down
#    to row 21 and keeps acting as the table's final row.
# ---------------------------------------------------------------------
$ws.Rows("19:20").Insert()

$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)
$ws.Range("B20:J20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Replace the previous account-statement rows with the new data set:
#    worker CC/1143379475 (LILIBETH ESTHER PEREZ BABILONIA) and worker
#    PPT/5614268 (JEFFERSON SMITH RIVERA CASTILLO), each with 3 periods.
# ---------------------------------------------------------------------
$ws.Range("B16").Value2 = "CC"
$ws.Range("C16").Value2 = "1143379475"
$ws.Range("D16").Value2 = "LILIBETH ESTHER PEREZ BABILONIA"
$ws.Range("E16").Value2 = "2507"
$ws.Range("F16").Value2 = 56940
$ws.Range("G16").Value2 = 1423500

$ws.Range("B17").Value2 = "CC"
$ws.Range("C17").Value2 = "1143379475"
$ws.Range("D17").Value2 = "LILIBETH ESTHER PEREZ BABILONIA"
$ws.Range("E17").Value2 = "2506"
$ws.Range("F17").Value2 = 56940
$ws.Range("G17").Value2 = 1423500

$ws.Range("B18").Value2 = "CC"
$ws.Range("C18").Value2 = "1143379475"
$ws.Range("D18").Value2 = "LILIBETH ESTHER PEREZ BABILONIA"
$ws.Range("E18").Value2 = "2505"
$ws.Range("F18").Value2 = 56940
$ws.Range("G18").Value2 = 1423500

$ws.Range("B19").Value2 = "PPT"
$ws.Range("C19").Value2 = "5614268"
$ws.Range("D19").Value2 = "JEFFERSON SMITH RIVERA CASTILLO"
$ws.Range("E19").Value2 = "2507"
$ws.Range("F19").Value2 = 56940
$ws.Range("G19").Value2 = 1423500

$ws.Range("B20").Value2 = "PPT"
$ws.Range("C20").Value2 = "5614268"
$ws.Range("D20").Value2 = "JEFFERSON SMITH RIVERA CASTILLO"
$ws.Range("E20").Value2 = "2506"
$ws.Range("F20").Value2 = 56940
$ws.Range("G20").Value2 = 1423500

$ws.Range("B21").Value2 = "PPT"
$ws.Range("C21").Value2 = "5614268"
$ws.Range("D21").Value2 = "JEFFERSON SMITH RIVERA CASTILLO"
$ws.Range("E21").Value2 = "2505"
$ws.Range("F21").Value2 = 56940
$ws.Range("G21").Value2 = 1423500
